# Update crypto price/volume snapshot (rows 2-51 of the cryptos sheet).
# "D" (Price) cells that look like plain numbers are written with a
# leading apostrophe so Excel stores them as text (matching the original
# t="inlineStr" string cells) instead of silently parsing them into
# numeric values; the Style reset afterwards clears the quote-prefix
# formatting flag that the apostrophe trick leaves behind so the cell's
# style index is left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.946.82"
$ws.Range("D3").Value = "1.638.05"
$ws.Range("E3").Value = "  -0.04%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'214.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("E8").Value = "  -0.32%  "
$ws.Range("D9").Value = "'0.0636"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.42%  "
$ws.Range("D10").Value = "'19.49"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.82%  "
$ws.Range("D11").Value = "'0.0795"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.00%  "
$ws.Range("D12").Value = "'4.24"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.34%  "
$ws.Range("D13").Value = "1.632.34"
$ws.Range("E13").Value = "  -0.58%  "
$ws.Range("D14").Value = "'0.541"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.70%  "
$ws.Range("D15").Value = "'63.19"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.82%  "
$ws.Range("E16").Value = "  -0.47%  "
$ws.Range("D17").Value = "25.981.24"
$ws.Range("E17").Value = "  +0.10%  "
$ws.Range("E18").Value = "  -0.01%  "
$ws.Range("D19").Value = "'194.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("E20").Value = "  -1.01%  "
$ws.Range("D21").Value = "'9.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.86%  "
$ws.Range("E22").Value = "  -1.86%  "
$ws.Range("D23").Value = "'0.132"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.78%  "
$ws.Range("D24").Value = "'143.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.30%  "
$ws.Range("E25").Value = "  -0.15%  "
$ws.Range("E26").Value = "  -0.65%  "
$ws.Range("E27").Value = "  +0.28%  "
$ws.Range("E28").Value = "  -0.25%  "
$ws.Range("E29").Value = "  -0.23%  "
$ws.Range("E30").Value = "  -1.60%  "
$ws.Range("E31").Value = "  -0.88%  "
$ws.Range("E32").Value = "  +0.42%  "
$ws.Range("E33").Value = "  -0.86%  "
$ws.Range("E34").Value = "  +0.65%  "
$ws.Range("D35").Value = "'0.900"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.40%  "
$ws.Range("D36").Value = "1.129.56"
$ws.Range("E36").Value = "  -0.79%  "
$ws.Range("E37").Value = "  -1.42%  "
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("E39").Value = "  -0.72%  "
$ws.Range("E40").Value = "  -0.88%  "
$ws.Range("D41").Value = "'5.41"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.27%  "
$ws.Range("D42").Value = "'0.792"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.99%  "
$ws.Range("E43").Value = "  -0.16%  "
$ws.Range("D44").Value = "'56.29"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.48%  "
$ws.Range("E45").Value = "  +2.83%  "
$ws.Range("E46").Value = "  -1.59%  "
$ws.Range("D47").Value = "'7.76"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.76%  "
$ws.Range("D48").Value = "'0.412"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.53%  "
$ws.Range("D49").Value = "'1.00"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.04%  "
$ws.Range("E50").Value = "  -2.04%  "
$ws.Range("D51").Value = "'5.49"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.74%  "
